$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ben@gmail.com"
$ws.Range("D3").Value = "gwen@gmail.com"
$ws.Range("I2").Value = "jarvis@gmail.com"
$ws.Range("I3").Value = "tony@gmail.com"

$ws.Range("A2").Value = "Ben"
$ws.Range("B2").Value = "Jarvis"
$ws.Range("A3").Value = "Gwen"
$ws.Range("B3").Value = "Tony"

$ws.Range("A4").Select()
